# Insert a new data row at row 211, shifting existing rows 211:319 down to 212:320.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with its values.
$ws.Range("A211").Value2 = 6
$ws.Range("B211").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C211").Value2 = "Metropolitana"
$ws.Range("D211").Value2 = 44529
$ws.Range("E211").Value2 = 13
$ws.Range("F211").Value2 = 100112043
$ws.Range("G211").Value2 = "Pepino ensalada"
$ws.Range("H211").Value2 = "Sin especificar"
$ws.Range("I211").Value2 = "Primera"
$ws.Range("J211").Value2 = 2200
$ws.Range("K211").Value2 = 5000
$ws.Range("L211").Value2 = 6000
$ws.Range("M211").Value2 = 5455
$ws.Range("N211").Value2 = "`$/caja 60 unidades"
$ws.Range("O211").Value2 = "Región de Arica y Parinacota"
$ws.Range("P211").Value2 = 91
$ws.Range("Q211").Value2 = 60
$ws.Range("R211").Value2 = "Hortaliza"
